$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the
# existing header row (B1:E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the time_taken values for each data row.
$times = @(
    "2021-10-05 13:41:12.529080",
    "2021-10-05 13:41:12.529090",
    "2021-10-05 13:41:12.529094",
    "2021-10-05 13:41:12.529096",
    "2021-10-05 13:41:12.529099",
    "2021-10-05 13:41:12.529102",
    "2021-10-05 13:41:12.529104",
    "2021-10-05 13:41:12.529107",
    "2021-10-05 13:41:12.529110",
    "2021-10-05 13:41:12.529112",
    "2021-10-05 13:41:12.529115",
    "2021-10-05 13:41:12.529118"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
